# BitcoinPriceTracker - SendEmail run without ForEach:
# shift current "Prices" (col B) into "Old Prices" (col D), then write freshly
# fetched Prices (col B) and Euro (col C) values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (CoinGecko) stores its price as text, so move it with a
#     values-only paste: this keeps it as a shared-string cell and keeps the
#     destination cell's (lack of) formatting intact, instead of coercing the
#     numeric-looking text back into a number or carrying over B6's style.
$ws.Cells.Item(6, 2).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)

# --- Rows 2-5: move the old numeric "Prices" into "Old Prices"
$ws.Cells.Item(2, 4).Value = 46371.15
$ws.Cells.Item(3, 4).Value = 46458.25
$ws.Cells.Item(4, 4).Value = 46457.43
$ws.Cells.Item(5, 4).Value = 46386.05

# --- Rows 2-5: newly fetched "Prices"
$ws.Cells.Item(2, 2).Value = 46263.47
$ws.Cells.Item(3, 2).Value = 46307.8
$ws.Cells.Item(4, 2).Value = 46264.35
$ws.Cells.Item(5, 2).Value = 46217.61

# --- Row 6: newly fetched "Prices" text value ("46238.19" + trailing NBSP,
#     matching the formatting already used for this column, e.g. cell D6's
#     original "46487.06\u00A0"). Build it in a scratch cell via a formula
#     (so it is treated as text, not a number), then paste its value into B6
#     so it becomes a brand new shared string while keeping B6's existing
#     style (s="2").
$ws.Cells.Item(1, 6).Formula = "=""46238.19" + [char]0x00A0 + """"
$ws.Cells.Item(1, 6).Copy()
$ws.Cells.Item(6, 2).PasteSpecial(-4163)
$ws.Cells.Item(1, 6).Clear()

# --- Newly fetched "Euro" values
$ws.Cells.Item(2, 3).Value = 40946.31
$ws.Cells.Item(3, 3).Value = 40985.55
$ws.Cells.Item(4, 3).Value = 40947.09
$ws.Cells.Item(5, 3).Value = 40905.72
$ws.Cells.Item(6, 3).Value = 40923.94
